$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Bump the revision number: "Version 1" -> "Version 2"
#    The document already stores "Version " and "1" as two separate
#    runs; only the second run's text changes, so only that run is
#    touched (leaving the "Version " run untouched, matching the
#    target markup). A harmless Bold on/off toggle keeps the engine
#    from silently re-merging the two runs back together.
# ---------------------------------------------------------------------
$r = $d.Content
$foundVersion = $r.Find.Execute("Version ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundVersion) {
    $numStart = $r.End
    $rNum = $d.Range($numStart, $numStart + 1)
    if ($rNum.Text -eq "1") {
        $rNum.Text = "2"
        $rNum.Font.Bold = 1
        $rNum.Font.Bold = 0
    }
}

# ---------------------------------------------------------------------
# 2. Update the generated date/time stamp text
# ---------------------------------------------------------------------
$r = $d.Content
[void]$r.Find.Execute("3/13/23 2:00 PM", $true, $false, $false, $false, $false, $true, 1, $false, "4/22/24 10:06 AM", 2)

# ---------------------------------------------------------------------
# 3. Shorten the first reference entry:
#    "Autonomous Vehicle Cybersecurity Development Lifecycle (AVCDL primary document)"
#    becomes
#    "AVCDL (primary document)"
#    The bold run that used to read "Autonomous Vehicle Cybersecurity
#    Development Lifecycle " is replaced by "AVCDL" followed by a
#    separate bold run containing just a space (matching how Word
#    splits a run when only part of it is replaced), and the trailing
#    parenthetical run drops the leading "AVCDL ".
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Autonomous Vehicle Cybersecurity Development Lifecycle", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Text = "AVCDL"
    $afterEnd = $r.End
    $rSpace = $d.Range($afterEnd, $afterEnd + 1)
    # Force the trailing space to live in its own run (rather than being
    # re-merged into the "AVCDL" run) while keeping it bold, matching the
    # target markup exactly.
    $rSpace.Font.Bold = 0
    $rSpace.Font.Bold = 1
}

$r = $d.Content
[void]$r.Find.Execute("(AVCDL primary document)", $true, $false, $false, $false, $false, $true, 1, $false, "(primary document)", 2)
